$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    5   = "18_hazards_to_humans_and_domestic_animals"
    6   = "ppe"
    7   = "ppe"
    8   = "ppe"
    12  = "env warning - water"
    13  = "application instructions || env warning - species"
    19  = "application instructions"
    20  = "application instructions"
    21  = "application instructions"
    22  = "135_product_information"
    23  = "use restrictions"
    25  = "application instructions"
    26  = "mixing"
    27  = "mixing"
    34  = "use restrictions"
    39  = "off target movement"
    40  = "off target movement"
    41  = "use restrictions || off target movement"
    42  = "application instructions"
    43  = "application instructions"
    44  = "application instructions"
    45  = "application instructions"
    46  = "application instructions"
    47  = "mixing"
    48  = "mixing"
    50  = "mixing"
    51  = "mixing"
    52  = "application instructions"
    53  = "application instructions"
    54  = "application instructions"
    55  = "irrigation || application instructions || chemigation"
    56  = "safety procedures"
    60  = "irrigation"
    61  = "irrigation"
    62  = "irrigation"
    64  = "irrigation"
    332 = "154_pesticide_storage"
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
